# "formulas" sheet: the used range grows from A1:D61 to A1:O61.
# For every data row (2-61) the existing "bundle" column (C, values
# vpc001..vpc060) gets mirrored into new column O, and the columns in
# between (E:N) are filled in with the same alternating row style
# (blank/empty cells) so the whole A1:O61 block is uniformly formatted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formulas")

$firstRow = 2
$lastRow = 61

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcCell = $ws.Range("C" + $r)

    # Copy column C's formatting (and, for O, its value) using the
    # clipboard so the existing style record is reused instead of a new
    # one being appended to the stylesheet.
    $srcCell.Copy()

    # E:N -> same row style, left blank/empty.
    $ws.Range("E" + $r + ":N" + $r).PasteSpecial(-4122)

    # O -> same row style, value copied from column C (bundle).
    $destCell = $ws.Range("O" + $r)
    $destCell.PasteSpecial(-4122)
    $destCell.Value = $srcCell.Value()
}

$excel.CutCopyMode = 0
